$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update column B (Read short-term) -----------------------
# Rule observed in the source data: any score of 11 or higher gets bumped
# up by one point (re-scoring adjustment); everything else is untouched.
for ($r = 2; $r -le 42; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    if ($b -ge 11) {
        $newB = $b + 1
    } else {
        $newB = $b
    }
    $ws.Cells.Item($r, 2).Value = $newB
}

# --- Step 2: add column P with the re-scoring formula -----------------
# P2 is a standalone formula; P3:P42 share the same formula (Excel
# auto-creates a shared formula group when a formula is filled down a
# contiguous range).
$ws.Range("P2").Formula = "=IF(B2>=11,B2+1,B2)"
$ws.Range("P3:P42").Formula = "=IF(B3>=11,B3+1,B3)"

# --- Step 3: add column Q with the literal re-scored values ------------
for ($r = 2; $r -le 42; $r++) {
    $newB = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 17).Value = $newB
}

# --- Step 4: selection cursor moved to E7 ------------------------------
$ws.Range("E7").Select()
